# "I finished the zombie and started adding audio"
#
# - I2 / I3 ("Figure out how to change direction" / "Reverse zombie
#   movement...") are marked done by giving them the same yellow highlight
#   already used elsewhere in the sheet (matches the other "Past"/done
#   cells' fill).
# - Two new tasks are logged in column I (the second "Tuesday" column):
#     I4 - "Start work on Zombie Animation"  (done-highlighted, like the
#          other single-line entries in that column)
#     I6 - "Make camera position dependent on player" (plain/open item,
#          wrapped like the rest of the backlog entries)
# - The active selection moves to I4, the newest/most relevant entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two already-existing Tuesday items as done: give them the
# yellow "done" fill while keeping their existing centered/wrapped look.
$ws.Range("I2").Interior.ColorIndex = 6
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("I2").WrapText = $true

$ws.Range("I3").Interior.ColorIndex = 6
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("I3").WrapText = $true

# New task: zombie animation work has started - highlight it done-style,
# centered, no wrap (single short line).
$ws.Range("I4").Value = "Start work on Zombie Animation"
$ws.Range("I4").Interior.ColorIndex = 6
$ws.Range("I4").HorizontalAlignment = -4108
$ws.Range("I4").WrapText = $false

# New task: camera-follows-player idea - plain (no fill) backlog entry,
# centered + wrapped like the other multi-line notes. Leave the fill
# untouched (cell default is already "no fill").
$ws.Range("I6").Value = "Make camera position dependent on player"
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("I6").WrapText = $true

# Widen column I so the new text is readable, mirroring the bestFit
# auto-sizing the rest of the sheet already uses (real Excel widens this
# column to just fit "Start work on Zombie Animation").
$ws.Columns.Item(9).ColumnWidth = 29.43

# Move the active selection to the newest entry.
$ws.Range("I4").Select() | Out-Null
